$d = $word.ActiveDocument

# This document repeats the same editorial markup pattern many times, so every
# deletion below is anchored on a short, verified-unique slice of plain text
# (found via Range.Find) rather than a bare "2", to make sure we only touch
# the three specific subscript "2" runs called out by the change - not any of
# the many other subscript/superscript "2" markers scattered through the
# index.

# --- Edit 1 ------------------------------------------------------------
# Paragraph "om. (1): 1/7c6<sub>2</sub><sub>β</sub> » [om.<sup>WH</sup><sub>2</sub>]"
# Remove the subscript "2" run that sits between "1/7c6" and the subscript "β".
$rng1 = $d.Content
$rng1.Find.Execute("1/7c62" + [char]946) | Out-Null
$two1 = $d.Range($rng1.End - 2, $rng1.End - 1)
if ($two1.Text -eq "2") {
    $two1.Delete()
}

# --- Edit 2 --------------------------------------------------------------
# Paragraph "ἄτυφος (1): 5/21a19<sub>β</sub><sup>WGH</sup><sub>2</sub> » невелан S<sub>2</sub>"
# Remove the subscript "2" run that sits between the superscript "WGH" and " »".
$rng2 = $d.Content
$rng2.Find.Execute("WGH2 " + [char]187) | Out-Null
$two2 = $d.Range($rng2.Start + 3, $rng2.Start + 4)
if ($two2.Text -eq "2") {
    $two2.Delete()
}

# --- Edit 3 --------------------------------------------------------------
# Paragraph "ἀκούω (1): 5/22b5<sub>2</sub><sub>β</sub> » [ѹслꙑшат<sup>GH</sup><sub>2</sub>]"
# Remove the subscript "2" run that sits between "5/22b5" and the subscript "β".
$rng3 = $d.Content
$rng3.Find.Execute("5/22b52" + [char]946) | Out-Null
$two3 = $d.Range($rng3.End - 2, $rng3.End - 1)
if ($two3.Text -eq "2") {
    $two3.Delete()
}
